# Applies the diff changes to the 'Bill Summary' worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Long text blocks, stored as single-quoted here-strings (no interpolation) ----
$e8text = @'
P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .
'@
$e10text = @'
Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   
'@
$e11text = @'
Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .
'@
$e12text = @'
Double pole MCB(With B/C curve tripping Characteristics)
'@

# Helper: force a cell to be stored as Text, then assign the value.
# (Needed for cells whose value looks numeric, e.g. "4.0", "17", "2050.00",
#  so Excel does not silently convert them to a number.)
function Set-TextValue {
    param(
        [object]$Cell,
        [string]$Text
    )
    $Cell.NumberFormat = "@"
    $Cell.Value2 = $Text
}

# ---------------------------------------------------------------------
# Row 8
Set-TextValue $ws.Range("A8") "Each"
$ws.Range("C8").Value2 = 41
Set-TextValue $ws.Range("D8") "4.0"
Set-TextValue $ws.Range("E8") $e8text
$ws.Range("F8").Value2 = 50
Set-TextValue $ws.Range("G8") "2050.00"

# Row 9
Set-TextValue $ws.Range("A9") "R. mtr."
$ws.Range("C9").Value2 = 84
Set-TextValue $ws.Range("D9") "17"
Set-TextValue $ws.Range("E9") "25 mm"
$ws.Range("F9").Value2 = 56
Set-TextValue $ws.Range("G9") "4704.00"

# Row 10
Set-TextValue $ws.Range("A10") "Set"
$ws.Range("C10").Value2 = 58
Set-TextValue $ws.Range("D10") "13.0"
Set-TextValue $ws.Range("E10") $e10text
$ws.Range("F10").Value2 = 5733
Set-TextValue $ws.Range("G10") "332514.00"

# Row 11
$ws.Range("C11").Value2 = 2
Set-TextValue $ws.Range("D11") "15.0"
Set-TextValue $ws.Range("E11") $e11text

# Row 12
$ws.Range("C12").Value2 = 22
Set-TextValue $ws.Range("D12") "31"
Set-TextValue $ws.Range("E12") $e12text

# Row 14 (Grand Total)
Set-TextValue $ws.Range("G14") "339268.00"
Set-TextValue $ws.Range("H14") "339268.00"

# Row 16 (Net Payable Amount)
Set-TextValue $ws.Range("G16") "339268.00"
Set-TextValue $ws.Range("H16") "339268.00"
